# Daily attendance processing - 2026-01-17 15:58:46
#
# The "Recorded By" column (G) lists the users/systems that contributed to
# recording each attendance session, as a comma-separated string. This pass
# re-syncs the display order with the backend's recency-first ordering: the
# most recently-acting contributor (previously listed last) is moved to the
# front, while the relative order of the remaining contributors is kept
# intact (i.e. the list is rotated right by one position).
#
# Rows whose "Recorded By" only has a single contributor are left untouched,
# since there is nothing to reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = [string]$cell.Text

    if ($text -eq '') { continue }

    $parts = $text -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1) {
        $rotated = @($parts[$parts.Length - 1]) + $parts[0..($parts.Length - 2)]
        $cell.Value = [string]::Join(', ', $rotated)
    }
}
